$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("match_number"), shifting the
# existing result/toss/match_link columns one to the right.
$ws.Columns("B:B").Insert()

# Populate the new column's header and first data row.
$ws.Range("B1").Value = "match_number"
$ws.Range("B2").Value = "Match 1"

# The new header cell (B1) should look like the other header cells
# (bold font + border, via the shared header style) - copy that
# formatting from the neighboring header cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats

# The new data cell (B2) should be plain/unstyled like its siblings.
$ws.Range("B2").ClearFormats()

$excel.CutCopyMode = $false
